# "sistemate le procedure invasive" - recompute/update the cost & effectiveness
# statistics (columns C:H) for each row of the results table after fixing the
# handling of invasive procedures. Row labels (column A) and infection counts
# (column B) are unchanged; only the derived statistics change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new C,D,E,F,G,H values (Eff_cost, c_low, c_high, Eff_ln, ln_low, ln_high)
$updates = @{
    2  = @(5807.99333333333, 4855.67000424764, 6760.31666241903, 53.9760435303685, 43.8447940458174, 64.8208552734572)
    3  = @(6644.35251798561, 3130.31491519274, 10158.3901207785, 44.6564839335076, 19.8770485409426, 74.5580042109418)
    4  = @(6970.97318007663, 4214.90984201603, 9727.03651813723, 58.2167158235049, 33.0964458856711, 88.0781188363096)
    5  = @(6619.54285714286, -417.086371763283, 13656.172086049, 56.1391912395674, 3.55265678560925, 135.430435082127)
    6  = @(5020.51176470588, 3295.51733422747, 6745.50619518429, 44.7355385388438, 29.9022197539786, 61.262649366602)
    7  = @(3619.74216027875, 1140.7151108601, 6098.76920969739, 19.8694162093934, 2.40293877720188, 40.315083863391)
    8  = @(6097.0652173913, 3063.92640999149, 9130.20402479112, 35.5654053995403, 13.3427839659655, 62.1451185340588)
    9  = @(6088.94926004228, 4492.15432076563, 7685.74419931894, 54.0900366192973, 38.3054790735635, 71.6760575530584)
    10 = @(7636.37384898711, 5876.46864893458, 9396.27904903964, 60.4975508251942, 43.8285646268481, 79.098386246965)
    11 = @(9950.73280943026, 7921.78186846451, 11979.683750396, 82.9410667004864, 62.3965904229516, 106.08458464767)
    12 = @(4480.01323529412, 3253.60221519012, 5706.42425539812, 49.7753475721206, 36.6815663265411, 64.1234830946879)
    13 = @(9557.82317073171, 7182.99580392167, 11932.6505375417, 90.8842070539401, 64.9825415136985, 120.852340910181)
    14 = @(5380.85972850679, 2663.67094096935, 8098.04851604422, 36.7634470939795, 17.02188460142, 59.8354062125637)
}

$columns = @("C", "D", "E", "F", "G", "H")

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
